$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log Friday's hours
$ws.Range("B6").Value = 5

# Highlight the Total row using Excel's built-in "Good" cell style, then
# frame it with a thin grey border
$totalRow = $ws.Range("A9:J9")
$totalRow.Style = "Good"
$totalRow.Borders.Color = 8355711
$totalRow.Borders.Weight = 2
$totalRow.Borders.LineStyle = 1

# Update the active selection
$ws.Range("B7").Select()
